$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,7).Value = 31.29437
$ws.Cells.Item(2,8).Value = 93.88310999999999
$ws.Cells.Item(2,9).Value = 0.2388439922596655
$ws.Cells.Item(2,10).Value = 0.2388439922596655
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,13).Value = 23.463916
$ws.Cells.Item(2,14).Value = 70.39174800000001
$ws.Cells.Item(2,15).Value = 0.1350973942042957
$ws.Cells.Item(2,16).Value = 0.1350973942042957
$ws.Cells.Item(2,17).Value = 734.28846895292
$ws.Cells.Item(2,18).Value = 6608.59622057628
$ws.Cells.Item(2,19).Value = 0.03226720097563178
$ws.Cells.Item(2,20).Value = 0.03226720097563178
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,7).Value = 31.29437
$ws.Cells.Item(3,8).Value = 93.88310999999999
$ws.Cells.Item(3,9).Value = 0.2388439922596655
$ws.Cells.Item(3,10).Value = 0.2388439922596655
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,13).Value = 87.90742766666666
$ws.Cells.Item(3,14).Value = 263.722283
$ws.Cells.Item(3,15).Value = 0.5061416180048238
$ws.Cells.Item(3,16).Value = 0.5061416180048238
$ws.Cells.Item(3,17).Value = 2751.007567148903
$ws.Cells.Item(3,18).Value = 24759.06810434013
$ws.Cells.Item(3,19).Value = 0.1208888846930387
$ws.Cells.Item(3,20).Value = 0.1208888846930387
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,7).Value = 31.29437
$ws.Cells.Item(4,8).Value = 93.88310999999999
$ws.Cells.Item(4,9).Value = 0.2388439922596655
$ws.Cells.Item(4,10).Value = 0.2388439922596655
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.06395099999999999
$ws.Cells.Item(4,14).Value = 0.191853
$ws.Cells.Item(4,15).Value = 0.0003682085060634768
$ws.Cells.Item(4,16).Value = 0.0003682085060634769
$ws.Cells.Item(4,17).Value = 2.001306255869999
$ws.Cells.Item(4,18).Value = 18.01175630283
$ws.Cells.Item(4,19).Value = 0.00008794438957216807
$ws.Cells.Item(4,20).Value = 0.00008794438957216809
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,7).Value = 31.29437
$ws.Cells.Item(5,8).Value = 93.88310999999999
$ws.Cells.Item(5,9).Value = 0.2388439922596655
$ws.Cells.Item(5,10).Value = 0.2388439922596655
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,13).Value = 62.24619000000001
$ws.Cells.Item(5,14).Value = 186.73857
$ws.Cells.Item(5,15).Value = 0.3583927792848171
$ws.Cells.Item(5,16).Value = 0.3583927792848171
$ws.Cells.Item(5,17).Value = 1947.9553009503
$ws.Cells.Item(5,18).Value = 17531.5977085527
$ws.Cells.Item(5,19).Value = 0.08559996220142287
$ws.Cells.Item(5,20).Value = 0.08559996220142287
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,7).Value = 33.86972933333333
$ws.Cells.Item(6,8).Value = 101.609188
$ws.Cells.Item(6,9).Value = 0.2584995758255442
$ws.Cells.Item(6,10).Value = 0.2584995758255442
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,13).Value = 23.463916
$ws.Cells.Item(6,14).Value = 70.39174800000001
$ws.Cells.Item(6,15).Value = 0.1350973942042957
$ws.Cells.Item(6,16).Value = 0.1350973942042957
$ws.Cells.Item(6,17).Value = 794.7164840200693
$ws.Cells.Item(6,18).Value = 7152.448356180625
$ws.Cells.Item(6,19).Value = 0.03492261909694676
$ws.Cells.Item(6,20).Value = 0.03492261909694676
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,7).Value = 33.86972933333333
$ws.Cells.Item(7,8).Value = 101.609188
$ws.Cells.Item(7,9).Value = 0.2584995758255442
$ws.Cells.Item(7,10).Value = 0.2584995758255442
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,13).Value = 87.90742766666666
$ws.Cells.Item(7,14).Value = 263.722283
$ws.Cells.Item(7,15).Value = 0.5061416180048238
$ws.Cells.Item(7,16).Value = 0.5061416180048238
$ws.Cells.Item(7,17).Value = 2977.400781459578
$ws.Cells.Item(7,18).Value = 26796.60703313621
$ws.Cells.Item(7,19).Value = 0.1308373935619016
$ws.Cells.Item(7,20).Value = 0.1308373935619016
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,7).Value = 33.86972933333333
$ws.Cells.Item(8,8).Value = 101.609188
$ws.Cells.Item(8,9).Value = 0.2584995758255442
$ws.Cells.Item(8,10).Value = 0.2584995758255442
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.06395099999999999
$ws.Cells.Item(8,14).Value = 0.191853
$ws.Cells.Item(8,15).Value = 0.0003682085060634768
$ws.Cells.Item(8,16).Value = 0.0003682085060634769
$ws.Cells.Item(8,17).Value = 2.166003060596
$ws.Cells.Item(8,18).Value = 19.494027545364
$ws.Cells.Item(8,19).Value = 0.00009518174263276607
$ws.Cells.Item(8,20).Value = 0.00009518174263276608
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,7).Value = 33.86972933333333
$ws.Cells.Item(9,8).Value = 101.609188
$ws.Cells.Item(9,9).Value = 0.2584995758255442
$ws.Cells.Item(9,10).Value = 0.2584995758255442
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,13).Value = 62.24619000000001
$ws.Cells.Item(9,14).Value = 186.73857
$ws.Cells.Item(9,15).Value = 0.3583927792848171
$ws.Cells.Item(9,16).Value = 0.3583927792848171
$ws.Cells.Item(9,17).Value = 2108.26160733124
$ws.Cells.Item(9,18).Value = 18974.35446598116
$ws.Cells.Item(9,19).Value = 0.09264438142406309
$ws.Cells.Item(9,20).Value = 0.09264438142406309
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,7).Value = 10.670404
$ws.Cells.Item(10,8).Value = 32.011212
$ws.Cells.Item(10,9).Value = 0.08143835106389757
$ws.Cells.Item(10,10).Value = 0.08143835106389757
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,13).Value = 23.463916
$ws.Cells.Item(10,14).Value = 70.39174800000001
$ws.Cells.Item(10,15).Value = 0.1350973942042957
$ws.Cells.Item(10,16).Value = 0.1350973942042957
$ws.Cells.Item(10,17).Value = 250.369463142064
$ws.Cells.Item(10,18).Value = 2253.325168278576
$ws.Cells.Item(10,19).Value = 0.01100210901702719
$ws.Cells.Item(10,20).Value = 0.01100210901702719
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,7).Value = 10.670404
$ws.Cells.Item(11,8).Value = 32.011212
$ws.Cells.Item(11,9).Value = 0.08143835106389757
$ws.Cells.Item(11,10).Value = 0.08143835106389757
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,13).Value = 87.90742766666666
$ws.Cells.Item(11,14).Value = 263.722283
$ws.Cells.Item(11,15).Value = 0.5061416180048238
$ws.Cells.Item(11,16).Value = 0.5061416180048238
$ws.Cells.Item(11,17).Value = 938.0077678041106
$ws.Cells.Item(11,18).Value = 8442.069910236996
$ws.Cells.Item(11,19).Value = 0.04121933877512599
$ws.Cells.Item(11,20).Value = 0.04121933877512599
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,7).Value = 10.670404
$ws.Cells.Item(12,8).Value = 32.011212
$ws.Cells.Item(12,9).Value = 0.08143835106389757
$ws.Cells.Item(12,10).Value = 0.08143835106389757
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 0.06395099999999999
$ws.Cells.Item(12,14).Value = 0.191853
$ws.Cells.Item(12,15).Value = 0.0003682085060634768
$ws.Cells.Item(12,16).Value = 0.0003682085060634769
$ws.Cells.Item(12,17).Value = 0.6823830062039999
$ws.Cells.Item(12,18).Value = 6.141447055836
$ws.Cells.Item(12,19).Value = 0.00002998629358151069
$ws.Cells.Item(12,20).Value = 0.00002998629358151069
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,7).Value = 10.670404
$ws.Cells.Item(13,8).Value = 32.011212
$ws.Cells.Item(13,9).Value = 0.08143835106389757
$ws.Cells.Item(13,10).Value = 0.08143835106389757
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,13).Value = 62.24619000000001
$ws.Cells.Item(13,14).Value = 186.73857
$ws.Cells.Item(13,15).Value = 0.3583927792848171
$ws.Cells.Item(13,16).Value = 0.3583927792848171
$ws.Cells.Item(13,17).Value = 664.1919947607601
$ws.Cells.Item(13,18).Value = 5977.727952846841
$ws.Cells.Item(13,19).Value = 0.02918691697816289
$ws.Cells.Item(13,20).Value = 0.02918691697816289
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Mfge8"
$ws.Cells.Item(14,3).Value = "Pdgfrb"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 55.18980966666666
$ws.Cells.Item(14,8).Value = 165.569429
$ws.Cells.Item(14,9).Value = 0.4212180808508926
$ws.Cells.Item(14,10).Value = 0.4212180808508926
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 23.463916
$ws.Cells.Item(14,14).Value = 70.39174800000001
$ws.Cells.Item(14,15).Value = 0.1350973942042957
$ws.Cells.Item(14,16).Value = 0.1350973942042957
$ws.Cells.Item(14,17).Value = 1294.969058074655
$ws.Cells.Item(14,18).Value = 11654.72152267189
$ws.Cells.Item(14,19).Value = 0.05690546511468993
$ws.Cells.Item(14,20).Value = 0.05690546511468993
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Mfge8"
$ws.Cells.Item(15,3).Value = "Pdgfrb"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 55.18980966666666
$ws.Cells.Item(15,8).Value = 165.569429
$ws.Cells.Item(15,9).Value = 0.4212180808508926
$ws.Cells.Item(15,10).Value = 0.4212180808508926
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 87.90742766666666
$ws.Cells.Item(15,14).Value = 263.722283
$ws.Cells.Item(15,15).Value = 0.5061416180048238
$ws.Cells.Item(15,16).Value = 0.5061416180048238
$ws.Cells.Item(15,17).Value = 4851.5942012096
$ws.Cells.Item(15,18).Value = 43664.34781088641
$ws.Cells.Item(15,19).Value = 0.2131960009747575
$ws.Cells.Item(15,20).Value = 0.2131960009747575
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Mfge8"
$ws.Cells.Item(16,3).Value = "Pdgfrb"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 55.18980966666666
$ws.Cells.Item(16,8).Value = 165.569429
$ws.Cells.Item(16,9).Value = 0.4212180808508926
$ws.Cells.Item(16,10).Value = 0.4212180808508926
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 0.6666666666666666
$ws.Cells.Item(16,13).Value = 0.06395099999999999
$ws.Cells.Item(16,14).Value = 0.191853
$ws.Cells.Item(16,15).Value = 0.0003682085060634768
$ws.Cells.Item(16,16).Value = 0.0003682085060634769
$ws.Cells.Item(16,17).Value = 3.529443517992999
$ws.Cells.Item(16,18).Value = 31.764991661937
$ws.Cells.Item(16,19).Value = 0.000155096080277032
$ws.Cells.Item(16,20).Value = 0.000155096080277032
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Mfge8"
$ws.Cells.Item(17,3).Value = "Pdgfrb"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 55.18980966666666
$ws.Cells.Item(17,8).Value = 165.569429
$ws.Cells.Item(17,9).Value = 0.4212180808508926
$ws.Cells.Item(17,10).Value = 0.4212180808508926
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 62.24619000000001
$ws.Cells.Item(17,14).Value = 186.73857
$ws.Cells.Item(17,15).Value = 0.3583927792848171
$ws.Cells.Item(17,16).Value = 0.3583927792848171
$ws.Cells.Item(17,17).Value = 3435.35537857517
$ws.Cells.Item(17,18).Value = 30918.19840717653
$ws.Cells.Item(17,19).Value = 0.1509615186811682
$ws.Cells.Item(17,20).Value = 0.1509615186811682
